# "fetch testdata from excel"
# Rework the Stock test-data sheet:
#  - rename the sheet Sheet1 -> Stock
#  - rename the header cells to the underscored convention
#  - trim the data set down to the first two rows (TC_001/RCOM, TC_002/TATAMOTORS)
#  - clear out the now unused rows 4 & 5 (keep their formatting)
#  - leave the selection on A4
#  - stamp a standard confidentiality footer on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Stock"

# Header row
$ws.Range("A1").Value = "Test_Case_ID"
$ws.Range("B1").Value = "Stock_Name"

# Row 2
$ws.Range("A2").Value = "TC_001"
$ws.Range("B2").Value = "RCOM"

# Row 3
$ws.Range("A3").Value = "TC_002"
$ws.Range("B3").Value = "TATAMOTORS"

# Rows 4 & 5 no longer hold data - clear the values but keep the styling
$ws.Range("A4:B4").ClearContents()
$ws.Range("A5:B5").ClearContents()

# Leave the cursor on A4
$ws.Range("A4").Select()

# Footer: blank line then the Cisco Confidential legal line
$ws.PageSetup.RightFooter = "`r&1#&""Calibri""&8&K000000 Cisco Confidential"
